$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("equilibrium_concentrations")
$ws.Range("A2").Value = 0.0000000992090565181129
$ws.Range("B2").Value = 0.000353446654489197
$ws.Range("C2").Value = 0.000000753337267674679
$ws.Range("D2").Value = 0.0000000000082431279606603
$ws.Range("E2").Value = 0.000000152562810448648
$ws.Range("A3").Value = 0.000105346173555645
$ws.Range("B3").Value = 0.000108191110679258
$ws.Range("C3").Value = 0.000244863809944392
$ws.Range("D3").Value = 0.00000284508101508837
$ws.Range("E3").Value = 0.000000000143675009480694
$ws.Range("A4").Value = 0.000358889149835025
$ws.Range("B4").Value = 0.0000391323448262809
$ws.Range("C4").Value = 0.00030172441801032
$ws.Range("D4").Value = 0.0000119432371646621
$ws.Range("E4").Value = 0.000000000042173502574039
$ws.Range("A5").Value = 0.000668426237493383
$ws.Range("B5").Value = 0.0000217002853363244
$ws.Range("C5").Value = 0.000311625645409508
$ws.Range("D5").Value = 0.0000229740704910569
$ws.Range("E5").Value = 0.00000000002264365405691
$ws.Range("A6").Value = 0.00132659238798426
$ws.Range("B6").Value = 0.0000105700849577871
$ws.Range("C6").Value = 0.000301252206660317
$ws.Range("D6").Value = 0.0000440777083830248
$ws.Range("E6").Value = 0.0000000000114093919288656
$ws.Range("A7").Value = 0.00263123903634047
$ws.Range("B7").Value = 0.00000478800252486154
$ws.Range("C7").Value = 0.000270663025538479
$ws.Range("D7").Value = 0.0000785489719366616
$ws.Range("E7").Value = 0.00000000000575227574360278
$ws.Range("A8").Value = 0.00661739630317633
$ws.Range("B8").Value = 0.00000143441341301523
$ws.Range("C8").Value = 0.000203927474063071
$ws.Range("D8").Value = 0.000148838112523943
$ws.Range("E8").Value = 0.00000000000228724588809908
$ws.Range("A9").Value = 0.0133759662765017
$ws.Range("B9").Value = 0.000000496827364661657
$ws.Range("C9").Value = 0.000142772620640885
$ws.Range("D9").Value = 0.000210630551994502
$ws.Range("E9").Value = 0.00000000000113155282926749

$ws = $wb.Worksheets.Item("absorbance_calc_abs_errors")
$ws.Range("C2").Value = 1.14941886695781
$ws.Range("D2").Value = 1.54289390498235
$ws.Range("E2").Value = 1.63781065083709
$ws.Range("F2").Value = 1.68384409257399
$ws.Range("G2").Value = 1.70338878985998
$ws.Range("H2").Value = 1.71146202865199
$ws.Range("I2").Value = 1.74346983026997
$ws.Range("J2").Value = 1.78869682511132
$ws.Range("C3").Value = 2.70361966945306
$ws.Range("D3").Value = 2.38168820371788
$ws.Range("E3").Value = 2.2546097327535
$ws.Range("F3").Value = 2.23501854459792
$ws.Range("G3").Value = 2.18355671706614
$ws.Range("H3").Value = 2.11159181386002
$ws.Range("I3").Value = 2.02114083957377
$ws.Range("J3").Value = 1.98878846322931
$ws.Range("C4").Value = 0.000418866957811126
$ws.Range("D4").Value = -0.000106095017651242
$ws.Range("E4").Value = -0.00418934916291391
$ws.Range("F4").Value = -0.00115590742601213
$ws.Range("G4").Value = 0.00238878985997548
$ws.Range("H4").Value = 0.00746202865199441
$ws.Range("I4").Value = -0.00653016973002862
$ws.Range("J4").Value = 0.00169682511131697
$ws.Range("C5").Value = 0.00161966945306347
$ws.Range("D5").Value = -0.00731179628211853
$ws.Range("E5").Value = 0.00260973275350151
$ws.Range("F5").Value = 0.00701854459791829
$ws.Range("G5").Value = -0.00244328293385809
$ws.Range("H5").Value = -0.00240818613997584
$ws.Range("I5").Value = 0.00114083957377176
$ws.Range("J5").Value = -0.000211536770685861

$ws = $wb.Worksheets.Item("absorbance_calc_rel_errors")
$ws.Range("C2").Value = 1.14941886695781
$ws.Range("D2").Value = 1.54289390498235
$ws.Range("E2").Value = 1.63781065083709
$ws.Range("F2").Value = 1.68384409257399
$ws.Range("G2").Value = 1.70338878985998
$ws.Range("H2").Value = 1.71146202865199
$ws.Range("I2").Value = 1.74346983026997
$ws.Range("J2").Value = 1.78869682511132
$ws.Range("C3").Value = 2.70361966945306
$ws.Range("D3").Value = 2.38168820371788
$ws.Range("E3").Value = 2.2546097327535
$ws.Range("F3").Value = 2.23501854459792
$ws.Range("G3").Value = 2.18355671706614
$ws.Range("H3").Value = 2.11159181386002
$ws.Range("I3").Value = 2.02114083957377
$ws.Range("J3").Value = 1.98878846322931
$ws.Range("C4").Value = 0.000364549136476176
$ws.Range("D4").Value = -0.0000687589226514855
$ws.Range("E4").Value = -0.00255136977034952
$ws.Range("F4").Value = -0.000685998472410757
$ws.Range("G4").Value = 0.00140434442091445
$ws.Range("H4").Value = 0.00437912479577137
$ws.Range("I4").Value = -0.00373152556001636
$ws.Range("J4").Value = 0.000949538394693323
$ws.Range("C5").Value = 0.000599433550356575
$ws.Range("D5").Value = -0.00306060957811575
$ws.Range("E5").Value = 0.00115885113388167
$ws.Range("F5").Value = 0.00315015466692921
$ws.Range("G5").Value = -0.00111769576114277
$ws.Range("H5").Value = -0.00113916089875868
$ws.Range("I5").Value = 0.000564772066223642
$ws.Range("J5").Value = -0.000106353328650508

$ws = $wb.Worksheets.Item("constants_evaluated")
$ws.Range("B4").Value = 4.3321142578125
$ws.Range("C4").Value = 0.0500355046429244
$ws.Range("B5").Value = 6.37466552734375
$ws.Range("C5").Value = 0.225001108736713

$ws = $wb.Worksheets.Item("correlation_matrix")
$ws.Range("B2").Value = 0.740739707949793
$ws.Range("A3").Value = 0.740739707949793

$ws = $wb.Worksheets.Item("adj_r_squared")
$ws.Range("A2").Value = 0.99984204419463

$ws = $wb.Worksheets.Item("mol_ext_coefficients_calc")
$ws.Range("C2").Value = 6.92335166639162
$ws.Range("D2").Value = 3241.77493862815
$ws.Range("E2").Value = 4810.10850140424
$ws.Range("F2").Value = 4784.3387312375
$ws.Range("C3").Value = 14.5869323319517
$ws.Range("D3").Value = 7635.87417879795
$ws.Range("E3").Value = 6297.31866792107
$ws.Range("F3").Value = 4229.18583923081
$ws.Range("C4").Value = 2.19081267294341
$ws.Range("D4").Value = 15.5464922272782
$ws.Range("E4").Value = 12.5561471931515
$ws.Range("F4").Value = 131.739463051048
$ws.Range("C5").Value = 2.1846099494123
$ws.Range("D5").Value = 15.5024763265327
$ws.Range("E5").Value = 12.5205976865155
$ws.Range("F5").Value = 131.366476589205
